$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "littleponyforever"
$ws.Range("B9").Value = "trixie image"

$ws.Range("B9").Select()
